$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: plain numeric id, two new strings
$ws.Range("A4").Value = 123
$ws.Range("B4").Value = "DAS"
$ws.Range("C4").Value = "AAA"

# Row 5: a date, a decimal number and a currency-formatted number
$ws.Range("A5").Value = 42715
$ws.Range("A5").NumberFormat = "DD/MM/YY"

$ws.Range("B5").Value = 23.05

$ws.Range("C5").Value = 50
$ws.Range("C5").NumberFormat = "[`$R`$-416]\ #,##0;\-[`$R`$-416]\ #,##0"

$null = $ws.Range("B9").Select()
